# Auto-committed on 2022/04/01 週五
#
# Adds a new "findYearMonthAll" lookup-function row to the DBS sheet:
#   A3 = findYearMonthAll
#   B3 = YearMonth =
# then leaves the selection on the DBS sheet at B4 (just below the newly
# typed row) and switches the active tab back to the DBD sheet.

$wb = $excel.ActiveWorkbook

$wsDBS = $wb.Worksheets.Item("DBS")
$wsDBS.Range("A3").Value = "findYearMonthAll"
$wsDBS.Range("B3").Value = "YearMonth ="
$wsDBS.Range("B4").Select()

$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBD.Activate()
